# Update DQ metrics and reports:
# Remove the "unambiguous_rdCase_rel_py_ipat" (column L) and
# "orphaCase_rel_py_ipat" (column M) metric columns from the DQ_Metrics sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DQ_Metrics")

$ws.Range("L:M").Delete()
